# Apply updated crypto price/volume figures to sheet1 (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.228.02'
$ws.Range("D3").Value = '3.066.74'
$ws.Range("E3").Value = '  +0.59%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '515.32'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.09%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '140.77'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.20%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  -1.78%  '
$ws.Range("E9").Value = '  +0.34%  '
$ws.Range("E10").Value = '  -1.52%  '
$ws.Range("E11").Value = '  -1.60%  '
$ws.Range("D12").Value = '3.597.31'
$ws.Range("E12").Value = '  +0.45%  '
$ws.Range("E13").Value = '  +2.57%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '25.51'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -5.03%  '
$ws.Range("E15").Value = '  -2.49%  '
$ws.Range("D16").Value = '57.321.58'
$ws.Range("E16").Value = '  +0.23%  '
$ws.Range("D17").Value = '3.070.60'
$ws.Range("E17").Value = '  +0.53%  '
$ws.Range("E18").Value = '  -1.03%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '13.03'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -2.83%  '
$ws.Range("E20").Value = '  -0.23%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '332.53'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -0.96%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '0.997'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -0.43%  '
$ws.Range("E23").Value = '  -1.38%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '65.59'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -0.06%  '
$ws.Range("E26").Value = '  -0.71%  '
$ws.Range("D27").Value = '0.0₃0907'
$ws.Range("E27").Value = '  -0.08%  '
$ws.Range("E28").Value = '  -6.08%  '
$ws.Range("E29").Value = '  -0.84%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '1.81'
$c.Style = "Normal"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '20.76'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +0.01%  '
$ws.Range("E32").Value = '  -2.91%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '154.97'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +1.13%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '27.34'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +7.94%  '
$ws.Range("E35").Value = '  -5.18%  '
$ws.Range("E36").Value = '  -2.04%  '
$ws.Range("E37").Value = '  -0.98%  '
$ws.Range("E38").Value = '  -0.45%  '
$ws.Range("D39").Value = '3.107.53'
$ws.Range("E39").Value = '  +0.60%  '
$ws.Range("E40").Value = '  -0.73%  '
$ws.Range("E41").Value = '  -0.56%  '
$ws.Range("E42").Value = '  +0.00%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.657'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -1.21%  '
$ws.Range("D44").Value = '2.270.75'
$ws.Range("E44").Value = '  +2.44%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.0260'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +6.72%  '
$ws.Range("E46").Value = '  -2.06%  '
$ws.Range("E47").Value = '  -2.65%  '
$ws.Range("E48").Value = '  -3.69%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '19.83'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -2.76%  '
$ws.Range("E50").Value = '  +0.26%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '248.19'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +5.08%  '
